# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.855.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "'1.790.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'310.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.5126"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "'0.3895"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.07822"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.76%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'41.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "'6.226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("D14").Value = "'20.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.782.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.224"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "'91.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "'0.00001075"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.75%  "
$ws.Range("D19").Value = "'0.06522"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'17.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").Value = "'5.919"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").Value = "'27.912.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "'11.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").Value = "'2.224"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "'160.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'20.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").Value = "'1.989.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").Value = "'2.354"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").Value = "'124.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'0.1074"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'1.039"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").Value = "'3.611"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").Value = "'5.493"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").Value = "'0.07056"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("D36").Value = "'0.02305"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Value = "'8.744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  -4.53%  "
$ws.Range("D39").Value = "'11.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "'4.994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("D41").Value = "'0.6090"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("D43").Value = "'1.146"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.307"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.79%  "
$ws.Range("D46").Value = "'0.5906"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "'3.696"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "'124.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").Value = "'1.203"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'1.911"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("D51").Value = "'0.06819"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.41%  "
